$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.199171543121338
$ws.Range("B1").Value = 4.499747753143311
$ws.Range("C1").Value = 3.362667322158813
$ws.Range("D1").Value = 0.8964079618453979
$ws.Range("E1").Value = 0.4716604351997375
